# "Generate Report for Handback"
#
# The localization-status report is refreshed after a handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (shown on the Overview sheet's zh-cn/de-de columns and on the per-locale
#    "Status" column).
#  - Each locale sheet's "Latest Target File" column now links to the source
#    markdown file, "Latest Handback File" records the generated xliff file,
#    and "Latest Handback DateTime" records when the handback happened
#    (handback times differ per locale).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFile = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9685f3c9a8c8df3d73d605eef853c75055fd5751/e2e/df11907f-0aea-4bd0-9a3d-fd2a40b6678b.md"

# --- Overview sheet: zh-cn / de-de status columns (rows 2 and 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn locale sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $targetFile
$wsZh.Range("I3").Value = $targetFile
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFile)

$wsZh.Range("J2").Value = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b.0b92c6b94efd9d28c696d9f7b188aa9549738b07.zh-cn.xlf"
$wsZh.Range("J3").Value = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b.0b92c6b94efd9d28c696d9f7b188aa9549738b07.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-09-02 19:10:57"
$wsZh.Range("K3").Value = "2016-09-02 19:10:57"

# --- de-de locale sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $targetFile
$wsDe.Range("I3").Value = $targetFile
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFile)

$wsDe.Range("J2").Value = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b.0b92c6b94efd9d28c696d9f7b188aa9549738b07.de-de.xlf"
$wsDe.Range("J3").Value = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b.0b92c6b94efd9d28c696d9f7b188aa9549738b07.de-de.xlf"

$wsDe.Range("K2").Value = "2016-09-02 19:11:13"
$wsDe.Range("K3").Value = "2016-09-02 19:11:13"

# --- Column widths widened to fit the new hyperlink / status text ---
# (ColumnWidth is in "characters"; the saved sheet <col width=.../> is
#  ColumnWidth + 5/6, so back that constant out to land on the desired
#  stored width.)
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527 - (5/6)
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527 - (5/6)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527 - (5/6)
$wsZh.Columns.Item(9).ColumnWidth = 40 - (5/6)
$wsZh.Columns.Item(10).ColumnWidth = 40 - (5/6)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527 - (5/6)
$wsDe.Columns.Item(9).ColumnWidth = 40 - (5/6)
$wsDe.Columns.Item(10).ColumnWidth = 40 - (5/6)

Write-Host "Handback report regenerated"
